$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27: num_customers 45 -> 46, retention_rate recalculated (C27/D27)
$ws.Range("C27").Value = 46
$ws.Range("E27").Value = 0.02042628774422735

# Row 36: num_customers 131 -> 132, retention_rate recalculated (C36/D36)
$ws.Range("C36").Value = 132
$ws.Range("E36").Value = 0.06839378238341969

# Row 37: num_customers 817 -> 822, cohort_size 817 -> 822 (retention_rate stays 1)
$ws.Range("C37").Value = 822
$ws.Range("D37").Value = 822
